# Auto-generated Excel COM-interop script to apply market-price refresh diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 58
$ws.Range("H58").Value = 1806
$ws.Range("I58").Value = 175.5
$ws.Range("J58").Value = 2711.8333
$ws.Range("K58").Value = 526.5
$ws.Range("L58").Value = 8135.499899999999
$ws.Range("M58").Value = -376.5
$ws.Range("N58").Value = -8435.499899999999

# Row 93
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

# Row 98
$ws.Range("H98").Value = 1165.8788
$ws.Range("I98").Value = 1006.25
$ws.Range("J98").Value = 2059.8
$ws.Range("K98").Value = 1006.25
$ws.Range("L98").Value = 2059.8
$ws.Range("M98").Value = 491.75
$ws.Range("N98").Value = -5055.8

# Row 122
$ws.Range("H122").Value = 1165.8788
$ws.Range("I122").Value = 1006.25
$ws.Range("J122").Value = 2059.8
$ws.Range("K122").Value = 3018.75
$ws.Range("L122").Value = 6179.400000000001
$ws.Range("M122").Value = -568.75
$ws.Range("N122").Value = -11079.4

# Row 125
$ws.Range("H125").Value = 3110
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 3110
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 27990
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -32910

# Row 137
$ws.Range("H137").Value = 1806.7
$ws.Range("I137").Value = 1742.2667
$ws.Range("J137").Value = 2000
$ws.Range("K137").Value = 5226.800099999999
$ws.Range("L137").Value = 6000
$ws.Range("M137").Value = -2676.800099999999
$ws.Range("N137").Value = -11100

# Row 138
$ws.Range("H138").Value = 1840.3611
$ws.Range("I138").Value = 1145.1904
$ws.Range("J138").Value = 2813.6
$ws.Range("K138").Value = 3435.5712
$ws.Range("L138").Value = 8440.799999999999
$ws.Range("M138").Value = 1704.4288
$ws.Range("N138").Value = -18720.8

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 3046.9656
$ws.Range("I61").Value = 2188.25
$ws.Range("J61").Value = 3374.0952
$ws.Range("K61").Value = 2188.25
$ws.Range("L61").Value = 3374.0952
$ws.Range("M61").Value = -1976.25
$ws.Range("N61").Value = -3798.0952

# Row 74
$ws.Range("H74").Value = 2382.9216
$ws.Range("I74").Value = 2739.6765
$ws.Range("K74").Value = 2739.6765
$ws.Range("M74").Value = -1865.6765

# Row 77
$ws.Range("H77").Value = 2382.9216
$ws.Range("I77").Value = 2739.6765
$ws.Range("K77").Value = 13698.3825
$ws.Range("M77").Value = -9330.3825

# Row 102
$ws.Range("H102").Value = 2058.1667
$ws.Range("I102").Value = 1969.8
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 1969.8
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = -347.8
$ws.Range("N102").Value = -5744

# Row 122
$ws.Range("H122").Value = 1665.6296
$ws.Range("I122").Value = 1552.9
$ws.Range("J122").Value = 1987.7142
$ws.Range("K122").Value = 4658.700000000001
$ws.Range("L122").Value = 5963.142599999999
$ws.Range("M122").Value = -2208.700000000001
$ws.Range("N122").Value = -10863.1426

# Row 132
$ws.Range("H132").Value = 5215.6
$ws.Range("I132").Value = 2081.7144
$ws.Range("J132").Value = 7304.857
$ws.Range("K132").Value = 6245.1432
$ws.Range("L132").Value = 21914.571
$ws.Range("M132").Value = -3715.1432
$ws.Range("N132").Value = -26974.571

# Row 136
$ws.Range("H136").Value = 3046.9656
$ws.Range("I136").Value = 2188.25
$ws.Range("J136").Value = 3374.0952
$ws.Range("K136").Value = 6564.75
$ws.Range("L136").Value = 10122.2856
$ws.Range("M136").Value = -4014.75
$ws.Range("N136").Value = -15222.2856

$ws = $wb.Worksheets.Item("BSM")
# Row 64
$ws.Range("H64").Value = 1370.8572
$ws.Range("J64").Value = 1233.3334
$ws.Range("L64").Value = 1233.3334
$ws.Range("N64").Value = -1683.3334

# Row 67
$ws.Range("H67").Value = 1370.8572
$ws.Range("J67").Value = 1233.3334
$ws.Range("L67").Value = 1233.3334
$ws.Range("N67").Value = -2793.3334

# Row 75
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("N75").ClearContents()

# Row 78
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("N78").ClearContents()

# Row 99
$ws.Range("H99").Value = 1926.3684
$ws.Range("I99").Value = 1342.1428
$ws.Range("J99").Value = 3562.2
$ws.Range("K99").Value = 1342.1428
$ws.Range("L99").Value = 3562.2
$ws.Range("M99").Value = 155.8571999999999
$ws.Range("N99").Value = -6558.2

# Row 101
$ws.Range("H101").Value = 20000
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

# Row 102
$ws.Range("H102").Value = 14508
$ws.Range("I102").Value = 1639
$ws.Range("J102").Value = 31666.666
$ws.Range("K102").Value = 1639
$ws.Range("L102").Value = 31666.666
$ws.Range("M102").Value = 1606
$ws.Range("N102").Value = -38156.666

# Row 134
$ws.Range("H134").Value = 5561.2324
$ws.Range("I134").Value = 2534.1765
$ws.Range("J134").Value = 7540.4614
$ws.Range("K134").Value = 7602.529500000001
$ws.Range("L134").Value = 22621.3842
$ws.Range("M134").Value = -5067.529500000001
$ws.Range("N134").Value = -27691.3842

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5954143
$ws.Range("I31").Value = 1311.1395
$ws.Range("J31").Value = 25644278
$ws.Range("K31").Value = 1311.1395
$ws.Range("L31").Value = 25644278
$ws.Range("M31").Value = -1016.1395
$ws.Range("N31").Value = -25644868

# Row 34
$ws.Range("H34").Value = 5954143
$ws.Range("I34").Value = 1311.1395
$ws.Range("J34").Value = 25644278
$ws.Range("K34").Value = 1311.1395
$ws.Range("L34").Value = 25644278
$ws.Range("M34").Value = -1109.1395
$ws.Range("N34").Value = -25644682

# Row 58
$ws.Range("H58").Value = 1434023.4
$ws.Range("I58").Value = 3022.76
$ws.Range("J58").Value = 5011525
$ws.Range("K58").Value = 3022.76
$ws.Range("L58").Value = 5011525
$ws.Range("M58").Value = -2819.76
$ws.Range("N58").Value = -5011931

# Row 132
$ws.Range("H132").Value = 3009.7778
$ws.Range("I132").Value = 1854
$ws.Range("J132").Value = 3934.4
$ws.Range("K132").Value = 5562
$ws.Range("L132").Value = 11803.2
$ws.Range("M132").Value = -3032
$ws.Range("N132").Value = -16863.2

# Row 134
$ws.Range("H134").Value = 2413.652
$ws.Range("I134").Value = 1193.4
$ws.Range("J134").Value = 3352.3076
$ws.Range("K134").Value = 3580.2
$ws.Range("L134").Value = 10056.9228
$ws.Range("M134").Value = -1045.2
$ws.Range("N134").Value = -15126.9228

# Row 136
$ws.Range("H136").Value = 1434023.4
$ws.Range("I136").Value = 3022.76
$ws.Range("J136").Value = 5011525
$ws.Range("K136").Value = 9068.280000000001
$ws.Range("L136").Value = 15034575
$ws.Range("M136").Value = -6518.280000000001
$ws.Range("N136").Value = -15039675

$ws = $wb.Worksheets.Item("CUL")
# Row 43
$ws.Range("H43").Value = 3710
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 3710
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 11130
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -11358

# Row 86
$ws.Range("H86").Value = 1562
$ws.Range("I86").Value = 465
$ws.Range("J86").Value = 2293.3333
$ws.Range("K86").Value = 1395
$ws.Range("L86").Value = 6879.999899999999
$ws.Range("M86").Value = -209
$ws.Range("N86").Value = -9251.999899999999

# Row 89
$ws.Range("H89").Value = 1562
$ws.Range("I89").Value = 465
$ws.Range("J89").Value = 2293.3333
$ws.Range("K89").Value = 4185
$ws.Range("L89").Value = 20639.9997
$ws.Range("M89").Value = 1743
$ws.Range("N89").Value = -32495.9997

# Row 98
$ws.Range("H98").Value = 2569.7
$ws.Range("I98").Value = 624.25
$ws.Range("J98").Value = 3866.6667
$ws.Range("K98").Value = 1872.75
$ws.Range("L98").Value = 11600.0001
$ws.Range("M98").Value = -374.75
$ws.Range("N98").Value = -14596.0001

# Row 131
$ws.Range("H131").Value = 385565.28
$ws.Range("I131").Value = 1429045.4
$ws.Range("J131").Value = 1125.2106
$ws.Range("K131").Value = 4287136.199999999
$ws.Range("L131").Value = 3375.6318
$ws.Range("M131").Value = -4282096.199999999
$ws.Range("N131").Value = -13455.6318

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 1439147.8
$ws.Range("I132").Value = 5210611
$ws.Range("J132").Value = 2399.9048
$ws.Range("K132").Value = 15631833
$ws.Range("L132").Value = 7199.714399999999
$ws.Range("M132").Value = -15629303
$ws.Range("N132").Value = -12259.7144

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 58835676
$ws.Range("I132").Value = 90926220
$ws.Range("J132").Value = 2999.6667
$ws.Range("K132").Value = 272778660
$ws.Range("L132").Value = 8999.000100000001
$ws.Range("M132").Value = -272776130
$ws.Range("N132").Value = -14059.0001

# Row 136
$ws.Range("H136").Value = 25002596
$ws.Range("I136").Value = 55556900
$ws.Range("K136").Value = 166670700
$ws.Range("M136").Value = -166668150

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 592
$ws.Range("J107").Value = 300
$ws.Range("L107").Value = 900
$ws.Range("N107").Value = -4740

# Row 132
$ws.Range("H132").Value = 3725.3684
$ws.Range("I132").Value = 5156.8
$ws.Range("J132").Value = 3214.1428
$ws.Range("K132").Value = 15470.4
$ws.Range("L132").Value = 9642.428400000001
$ws.Range("M132").Value = -12940.4
$ws.Range("N132").Value = -14702.4284

# Row 136
$ws.Range("H136").Value = 7578081.5
$ws.Range("I136").Value = 14707534
$ws.Range("J136").Value = 3037.5
$ws.Range("K136").Value = 44122602
$ws.Range("L136").Value = 9112.5
$ws.Range("M136").Value = -44120052
$ws.Range("N136").Value = -14212.5

Write-Host "Applied all market data updates."